$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices and 1h volume changes).
# Columns B/C/D/E hold coin name, link, price, and volume-change text respectively.
# D-column price values that look like plain numbers must be forced to Text format
# first, otherwise Excel auto-converts them to floating point numbers and corrupts
# the original text representation (e.g. "0.0990" -> 0.099).

$ws.Range("D2").Value = "43.775.46"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.338.42"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.667"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.72"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.66"
$ws.Range("E7").Value = "  -3.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  +6.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.19"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "31.86"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.17"
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "2.685.33"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "2.331.28"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "43.590.48"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.26"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.36"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +23.11%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.53"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.47"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.39"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.55"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.71"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.18"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  +10.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("E42").Value = "  +8.50%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.89"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.41"
$ws.Range("E45").Value = "  +14.56%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.61"
$ws.Range("E46").Value = "  -3.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.67"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.44"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.06"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  -2.26%  "
